$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @(
    "Thu Nov 07 16:05:16 EST 2024",
    "Thu Nov 07 16:05:27 EST 2024",
    "Thu Nov 07 16:05:39 EST 2024",
    "Thu Nov 07 16:05:50 EST 2024",
    "Thu Nov 07 16:06:01 EST 2024",
    "Thu Nov 07 16:06:12 EST 2024",
    "Thu Nov 07 16:06:23 EST 2024",
    "Thu Nov 07 16:06:34 EST 2024",
    "Thu Nov 07 16:06:45 EST 2024",
    "Thu Nov 07 16:06:58 EST 2024",
    "Thu Nov 07 16:07:09 EST 2024",
    "Thu Nov 07 16:07:20 EST 2024",
    "Thu Nov 07 16:07:31 EST 2024",
    "Thu Nov 07 16:07:42 EST 2024",
    "Thu Nov 07 16:07:54 EST 2024",
    "Thu Nov 07 16:08:06 EST 2024",
    "Thu Nov 07 16:08:18 EST 2024",
    "Thu Nov 07 16:08:29 EST 2024",
    "Thu Nov 07 16:08:40 EST 2024",
    "Thu Nov 07 16:08:51 EST 2024",
    "Thu Nov 07 16:09:02 EST 2024",
    "Thu Nov 07 16:09:13 EST 2024",
    "Thu Nov 07 16:09:25 EST 2024",
    "Thu Nov 07 16:09:37 EST 2024",
    "Thu Nov 07 16:09:48 EST 2024",
    "Thu Nov 07 16:10:00 EST 2024",
    "Thu Nov 07 16:10:12 EST 2024",
    "Thu Nov 07 16:10:23 EST 2024",
    "Thu Nov 07 16:10:34 EST 2024"
)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newDates[$i]
}
